$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for updated rows.
# D-column values that look numeric need NumberFormat "@" first so they
# stay literal text (matching the original inlineStr cell type) instead of
# being auto-converted to a number by Excel.

$ws.Cells.Item(2, 4).Value = "46.373.08"
$ws.Cells.Item(2, 5).Value = "  -1.08%  "
$ws.Cells.Item(3, 4).Value = "2.461.96"
$ws.Cells.Item(3, 5).Value = "  +8.40%  "
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "1.00"
$ws.Cells.Item(4, 5).Value = "  +0.10%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "297.26"
$ws.Cells.Item(5, 5).Value = "  -0.89%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "96.18"
$ws.Cells.Item(6, 5).Value = "  -2.94%  "
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.576"
$ws.Cells.Item(7, 5).Value = "  +0.29%  "
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "1.00"
$ws.Cells.Item(8, 5).Value = "  +0.18%  "
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.517"
$ws.Cells.Item(9, 5).Value = "  +2.33%  "
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "35.30"
$ws.Cells.Item(10, 5).Value = "  +0.70%  "
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.0789"
$ws.Cells.Item(11, 5).Value = "  -0.81%  "
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "7.26"
$ws.Cells.Item(12, 5).Value = "  +3.26%  "
$ws.Cells.Item(13, 5).Value = "  +1.95%  "
$ws.Cells.Item(14, 4).Value = "2.844.01"
$ws.Cells.Item(14, 5).Value = "  +8.68%  "
$ws.Cells.Item(15, 4).Value = "2.478.84"
$ws.Cells.Item(15, 5).Value = "  +9.39%  "
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "0.857"
$ws.Cells.Item(16, 5).Value = "  +7.79%  "
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "14.28"
$ws.Cells.Item(17, 5).Value = "  +4.64%  "
$ws.Cells.Item(18, 4).Value = "46.389.45"
$ws.Cells.Item(18, 5).Value = "  -1.04%  "
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "12.77"
$ws.Cells.Item(19, 5).Value = "  +2.57%  "
$ws.Cells.Item(20, 4).Value = "0.0₃0950"
$ws.Cells.Item(20, 5).Value = "  -2.35%  "
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "6.32"
$ws.Cells.Item(21, 5).Value = "  +8.81%  "
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "67.99"
$ws.Cells.Item(22, 5).Value = "  +3.38%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "246.68"
$ws.Cells.Item(23, 5).Value = "  +0.49%  "
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "2.82"
$ws.Cells.Item(24, 5).Value = "  +1.14%  "
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "1.96"
$ws.Cells.Item(25, 5).Value = "  +5.48%  "
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "0.999"
$ws.Cells.Item(26, 5).Value = "  -0.15%  "
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "39.72"
$ws.Cells.Item(27, 5).Value = "  -4.08%  "
$ws.Cells.Item(28, 5).Value = "  +0.77%  "
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "9.87"
$ws.Cells.Item(29, 5).Value = "  +3.27%  "
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "21.73"
$ws.Cells.Item(30, 5).Value = "  +8.42%  "
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "3.84"
$ws.Cells.Item(31, 5).Value = "  +14.74%  "
$ws.Cells.Item(32, 5).Value = "  -1.72%  "
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "5.59"
$ws.Cells.Item(33, 5).Value = "  +5.13%  "
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "148.29"
$ws.Cells.Item(34, 5).Value = "  +2.05%  "
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "2.06"
$ws.Cells.Item(35, 5).Value = "  +23.80%  "
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "0.0774"
$ws.Cells.Item(36, 5).Value = "  +0.97%  "
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "0.116"
$ws.Cells.Item(37, 5).Value = "  +4.12%  "
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.116"
$ws.Cells.Item(38, 5).Value = "  +0.55%  "
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "15.28"
$ws.Cells.Item(39, 5).Value = "  -1.61%  "
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "3.95"
$ws.Cells.Item(40, 5).Value = "  +3.44%  "
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "0.0303"
$ws.Cells.Item(41, 5).Value = "  +2.81%  "
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "3.27"
$ws.Cells.Item(42, 5).Value = "  +6.27%  "
$ws.Cells.Item(43, 4).Value = "2.010.69"
$ws.Cells.Item(43, 5).Value = "  +13.17%  "
$ws.Cells.Item(44, 5).Value = "  -0.02%  "
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "92.47"
$ws.Cells.Item(45, 5).Value = "  -1.26%  "
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "8.62"
$ws.Cells.Item(48, 5).Value = "  +9.51%  "
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "103.17"
$ws.Cells.Item(49, 5).Value = "  +9.37%  "
$ws.Cells.Item(50, 4).Value = "2.710.17"
$ws.Cells.Item(50, 5).Value = "  +8.67%  "
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "0.187"
$ws.Cells.Item(51, 5).Value = "  +2.11%  "

# Rows 46 and 47: the Stacks / EnergySwap entries swapped rank position.
$ws.Cells.Item(46, 2).Value = "EnergySwap"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "16.60"
$ws.Cells.Item(46, 5).Value = "  +33.62%  "
$ws.Cells.Item(47, 2).Value = "Stacks"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "1.80"
$ws.Cells.Item(47, 5).Value = "  -4.21%  "
